# Feb 15th commit #4
# Update "Final Result" columns on each sheet with new timestamped Success values.

$wb = $excel.ActiveWorkbook

$wsSchoolSearch   = $wb.Worksheets.Item("School Search")
$wsProductSearch  = $wb.Worksheets.Item("Product Search")
$wsShoppingCart   = $wb.Worksheets.Item("Shopping Cart")
$wsCheckout       = $wb.Worksheets.Item("Checkout")
$wsPayment        = $wb.Worksheets.Item("Payment")

# School Search sheet ("Final Result" column C)
$wsSchoolSearch.Range("C2").Value = "Success - 2021/02/15 20:47:57"
$wsSchoolSearch.Range("C3").Value = "Success - 2021/02/15 20:48:00"

# Product Search sheet ("Final Result" column K)
$wsProductSearch.Range("K1").Value = "Success - 2021/02/15 20:49:27"

# Shopping Cart sheet ("Final Result" column G)
$wsShoppingCart.Range("G2").Value = "Success - 2021/02/15 20:49:30"
$wsShoppingCart.Range("G3").Value = "Success - 2021/02/15 20:49:30"
$wsShoppingCart.Range("G4").Value = "Success - 2021/02/15 20:49:30"

# Checkout sheet ("Final Result" column P)
$wsCheckout.Range("P2").Value = "Success - 2021/02/15 20:49:39"
$wsCheckout.Range("P3").Value = "Success - 2021/02/15 20:49:50"
$wsCheckout.Range("P4").Value = "Success - 2021/02/15 20:49:58"

# Payment sheet ("Final Result" column F)
$wsPayment.Range("F2").Value = "Success - 2021/02/15 20:50:14"
